# TC38_Canine_Filter_Breed-Poodle.xlsx -- "updated first 25 tc in icdc breed+diagnosis"
#
# The "startup" sheet drives the automation tool: column A = tab name,
# column B = the "query" run against the neo4j DB for that tab, column C =
# the "StatQuery" used to populate summary stats for that tab.
#
# This edit:
#   * Adds a `Cohort` column to the CasesTab query (B2).
#   * Replaces the shared StatQuery (C2:C4) with a new program/study/case/
#     sample/file counts query.
#   * Moves the old FilesTab query (which used to live in the StatQuery
#     slot's sibling) back into the FilesTab row's query cell (B4) -- i.e.
#     restores the original file-listing query for FilesTab, now that the
#     StatQuery column carries the new stats query instead.
#   * SamplesTab's query (B3) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New CasesTab query (B2): appends the Cohort column -------------------
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Poodle']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
coalesce(co.cohort_description, '') AS `Cohort`
'@

# --- New shared StatQuery (C2, C3, C4) -------------------------------------
$statQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Poodle']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

# --- FilesTab query (B4): restored file-listing query ----------------------
$filesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Poodle']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
         coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
         coalesce(diag.disease_term,'') AS Diagnosis   
'@

# Trim the single trailing newline the here-string literal adds (the source
# text itself has no trailing newline).
$casesQuery = $casesQuery -replace "`r?`n$", ''
$statQuery  = $statQuery  -replace "`r?`n$", ''
$filesQuery = $filesQuery -replace "`r?`n$", ''

# Write order matters for the position new entries take in the shared
# string table: C2 (stat query) first, then B2 (cases query) with the
# cohort line, then the FilesTab query back into B4.
$ws.Range("C2").Value = $statQuery
$ws.Range("B2").Value = $casesQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("B4").Value = $filesQuery
$ws.Range("C4").Value = $statQuery

# Row heights re-fit themselves to the new wrapped text.
$ws.Rows(2).RowHeight = 270
$ws.Rows(3).RowHeight = 225
$ws.Rows(4).RowHeight = 210

# The active selection moved to C2 in the saved workbook.
$ws.Range("C2").Select()
